# Add a new game-stat row for Rui Hachimura (inserted above the "promedios"
# average row, which shifts from row 3 down to row 4) and update that row's
# AVERAGE() formulas to include the newly inserted row. Because this is a
# plain row insertion, Excel automatically re-points every cross-sheet
# formula on the "final" summary sheet that referenced the old average row
# (Rui Hachimura!A3:I3) to the new one (Rui Hachimura!A4:I4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rui Hachimura")

# Push the existing "promedios" average row (row 3) down to row 4, leaving a
# fresh, blank row 3 for the new game's stats.
$ws.Rows.Item(3).Insert()

# New game stats for Rui Hachimura.
$ws.Range("A3").Value = 50
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 550
$ws.Range("D3").Value = 277
$ws.Range("E3").Value = 161
$ws.Range("F3").Value = 9.09
$ws.Range("G3").Value = 9.55
$ws.Range("H3").Value = 11.98
$ws.Range("I3").Value = -511

# Recompute the averages row (now row 4) over the two games (row 2 and the
# newly added row 3).
$ws.Range("A4").Formula = "=AVERAGE(A2:A3)"
$ws.Range("B4").Formula = "=AVERAGE(B2:B3)"
$ws.Range("C4").Formula = "=AVERAGE(C2:C3)"
$ws.Range("D4").Formula = "=AVERAGE(D2:D3)"
$ws.Range("E4").Formula = "=AVERAGE(E2:E3)"
$ws.Range("F4").Formula = "=AVERAGE(F2:F3)"
$ws.Range("G4").Formula = "=AVERAGE(G2:G3)"
$ws.Range("H4").Formula = "=AVERAGE(H2:H3)"
$ws.Range("I4").Formula = "=AVERAGE(I2:I3)"
